$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting so numeric-looking strings (prices, +/- changes,
# percentages) are stored as text, matching the original inline-string
# cell layout instead of being auto-converted to numbers by Excel.
$textRange = $ws.Range("B2:E8")
$textRange.NumberFormat = "@"

# New header for column E
$ws.Range("E1").Value = "Yahoo Estimated Return"

# Row 2 - AAPL
$ws.Range("B2").Value = "177.70"
$ws.Range("C2").Value = "+0.21"
$ws.Range("D2").Value = "+0.12%"
$ws.Range("E2").Value = "-2%"

# Row 3 - GOOG
$ws.Range("B3").Value = "138.71"
$ws.Range("C3").Value = "-0.02"
$ws.Range("D3").Value = "-0.02%"
$ws.Range("E3").Value = "7%"

# Row 4 - GOOGL
$ws.Range("B4").Value = "137.65"
$ws.Range("C4").Value = "+0.07"
$ws.Range("D4").Value = "+0.05%"
$ws.Range("E4").Value = "7%"

# Row 5 - AMZN
$ws.Range("B5").Value = "126.45"
$ws.Range("C5").Value = "-1.51"
$ws.Range("D5").Value = "-1.18%"
$ws.Range("E5").Value = "-16%"

# Row 6 - META
$ws.Range("B6").Value = "316.39"
$ws.Range("C6").Value = "+0.96"
$ws.Range("D6").Value = "+0.30%"
$ws.Range("E6").Value = "-5%"

# Row 7 - MSFT
$ws.Range("B7").Value = "327.56"
$ws.Range("C7").Value = "+0.30"
$ws.Range("D7").Value = "+0.09%"
$ws.Range("E7").Value = "-3%"

# Row 8 - NVDA
$ws.Range("B8").Value = "446.82"
$ws.Range("C8").Value = "-10.80"
$ws.Range("D8").Value = "-2.36%"
$ws.Range("E8").Value = "5%"
